# Applies the "Updated cryptos list" data refresh described in the commit
# message (GitHub Actions price-scrape run). For every affected row, the
# Price (D) and Volume(1h) (E) text is refreshed to the latest scrape, and
# three rows (31/32 and 38/39/40) have their Coin/Link/Price/Volume swapped
# to reflect the new ranking order.
#
# D/E hold plain text (e.g. "30.714.33", "  +2.33%  ") rather than numbers,
# so the whole D2:E51 block is forced to text format up front (matching how
# these "numeric-looking" strings, e.g. "80.00"/"96.00", were authored) and
# the style is restored back to Normal afterwards so no visible formatting
# changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$updates = @(
    @{ Row=2; D='30.714.33'; E='  +2.33%  ' },
    @{ Row=3; D='1.693.13'; E='  +3.53%  ' },
    @{ Row=4; D='0.997'; E='  -0.15%  ' },
    @{ Row=5; D='221.81'; E='  +3.06%  ' },
    @{ Row=6; D='0.524'; E='  +0.59%  ' },
    @{ Row=7; D='0.997'; E='  -0.18%  ' },
    @{ Row=8; D='31.17'; E='  +4.22%  ' },
    @{ Row=9; E='  +2.44%  ' },
    @{ Row=10; D='0.0629'; E='  +2.51%  ' },
    @{ Row=11; D='0.0902'; E='  -1.63%  ' },
    @{ Row=12; D='1.935.45'; E='  +3.55%  ' },
    @{ Row=13; D='10.83'; E='  +12.78%  ' },
    @{ Row=14; D='0.622'; E='  +7.88%  ' },
    @{ Row=15; D='1.689.17'; E='  +3.39%  ' },
    @{ Row=16; D='4.03'; E='  +3.35%  ' },
    @{ Row=17; D='30.701.03'; E='  +2.26%  ' },
    @{ Row=18; D='66.27'; E='  +2.11%  ' },
    @{ Row=19; D='248.70'; E='  -0.11%  ' },
    @{ Row=20; D='0.0₃0722'; E='  +1.94%  ' },
    @{ Row=21; E='  -0.18%  ' },
    @{ Row=22; D='4.31'; E='  +2.90%  ' },
    @{ Row=23; D='10.23'; E='  +5.50%  ' },
    @{ Row=24; D='2.18'; E='  +2.51%  ' },
    @{ Row=25; D='157.33'; E='  -1.57%  ' },
    @{ Row=26; D='15.99'; E='  +1.59%  ' },
    @{ Row=27; E='  +0.50%  ' },
    @{ Row=28; D='6.78'; E='  +1.99%  ' },
    @{ Row=29; D='0.997'; E='  -0.16%  ' },
    @{ Row=30; D='0.0502'; E='  +2.34%  ' },
    @{ Row=31; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='1.14'; E='  +1.40%  ' },
    @{ Row=32; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='3.50'; E='  +3.27%  ' },
    @{ Row=33; D='3.32'; E='  +3.14%  ' },
    @{ Row=34; D='1.515.54'; E='  +5.73%  ' },
    @{ Row=35; E='  +5.16%  ' },
    @{ Row=36; E='  -0.40%  ' },
    @{ Row=37; D='0.0181'; E='  +5.16%  ' },
    @{ Row=38; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='80.00'; E='  +8.04%  ' },
    @{ Row=39; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.588'; E='  +5.34%  ' },
    @{ Row=40; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.71'; E='  -5.64%  ' },
    @{ Row=41; E='  +1.36%  ' },
    @{ Row=42; D='0.858'; E='  +2.49%  ' },
    @{ Row=43; D='2.03'; E='  +1.63%  ' },
    @{ Row=44; D='0.0502'; E='  +0.71%  ' },
    @{ Row=45; E='  -1.58%  ' },
    @{ Row=46; D='0.997'; E='  -0.16%  ' },
    @{ Row=47; D='52.61'; E='  -5.42%  ' },
    @{ Row=48; D='1.827.51'; E='  +2.83%  ' },
    @{ Row=49; D='5.45'; E='  +0.35%  ' },
    @{ Row=50; D='96.00'; E='  +6.25%  ' },
    @{ Row=51; D='0.0₆0115'; E='  +5.43%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Cells.Item($row, 4).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($row, 5).Value = $u.E }
}

# Restore the default (unstyled) look for the text-formatted block.
$priceVolumeRange.Style = "Normal"

Write-Host "Applied $($updates.Count) row updates"
